{"js": "// Replace the 25 division-fact answers in the table with the new set,\n// matched on their unique current text.\nconst replacements = [\n  [\"26\u00f74=6, 2\", \"98\u00f77=14, 0\"],\n  [\"43\u00f74=10, 3\", \"97\u00f75=19, 2\"],\n  [\"36\u00f78=4, 4\", \"33\u00f72=16, 1\"],\n  [\"61\u00f74=15, 1\", \"75\u00f75=15, 0\"],\n  [\"22\u00f78=2, 6\", \"94\u00f77=13, 3\"],\n  [\"62\u00f72=31, 0\", \"15\u00f77=2, 1\"],\n  [\"31\u00f74=7, 3\", \"81\u00f76=13, 3\"],\n  [\"84\u00f79=9, 3\", \"87\u00f78=10, 7\"],\n  [\"61\u00f78=7, 5\", \"99\u00f77=14, 1\"],\n  [\"69\u00f79=7, 6\", \"16\u00f73=5, 1\"],\n  [\"38\u00f75=7, 3\", \"70\u00f77=10, 0\"],\n  [\"66\u00f77=9, 3\", \"85\u00f79=9, 4\"],\n  [\"99\u00f78=12, 3\", \"64\u00f74=16, 0\"],\n  [\"35\u00f79=3, 8\", \"39\u00f79=4, 3\"],\n  [\"14\u00f73=4, 2\", \"31\u00f74=7, 3\"],\n  [\"34\u00f73=11, 1\", \"52\u00f73=17, 1\"],\n  [\"39\u00f78=4, 7\", \"57\u00f76=9, 3\"],\n  [\"12\u00f73=4, 0\", \"94\u00f72=47, 0\"],\n  [\"60\u00f75=12, 0\", \"65\u00f73=21, 2\"],\n  [\"67\u00f76=11, 1\", \"51\u00f74=12, 3\"],\n  [\"15\u00f79=1, 6\", \"31\u00f77=4, 3\"],\n  [\"23\u00f73=7, 2\", \"82\u00f78=10, 2\"],\n  [\"55\u00f72=27, 1\", \"18\u00f74=4, 2\"],\n  [\"91\u00f72=45, 1\", \"79\u00f73=26, 1\"],\n  [\"20\u00f79=2, 2\", \"74\u00f73=24, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 division-fact answers in the table with the new set,\n# matched on their unique current text via Find/Replace on the whole story.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"26\u00f74=6, 2\", \"98\u00f77=14, 0\"),\n  @(\"43\u00f74=10, 3\", \"97\u00f75=19, 2\"),\n  @(\"36\u00f78=4, 4\", \"33\u00f72=16, 1\"),\n  @(\"61\u00f74=15, 1\", \"75\u00f75=15, 0\"),\n  @(\"22\u00f78=2, 6\", \"94\u00f77=13, 3\"),\n  @(\"62\u00f72=31, 0\", \"15\u00f77=2, 1\"),\n  @(\"31\u00f74=7, 3\", \"81\u00f76=13, 3\"),\n  @(\"84\u00f79=9, 3\", \"87\u00f78=10, 7\"),\n  @(\"61\u00f78=7, 5\", \"99\u00f77=14, 1\"),\n  @(\"69\u00f79=7, 6\", \"16\u00f73=5, 1\"),\n  @(\"38\u00f75=7, 3\", \"70\u00f77=10, 0\"),\n  @(\"66\u00f77=9, 3\", \"85\u00f79=9, 4\"),\n  @(\"99\u00f78=12, 3\", \"64\u00f74=16, 0\"),\n  @(\"35\u00f79=3, 8\", \"39\u00f79=4, 3\"),\n  @(\"14\u00f73=4, 2\", \"31\u00f74=7, 3\"),\n  @(\"34\u00f73=11, 1\", \"52\u00f73=17, 1\"),\n  @(\"39\u00f78=4, 7\", \"57\u00f76=9, 3\"),\n  @(\"12\u00f73=4, 0\", \"94\u00f72=47, 0\"),\n  @(\"60\u00f75=12, 0\", \"65\u00f73=21, 2\"),\n  @(\"67\u00f76=11, 1\", \"51\u00f74=12, 3\"),\n  @(\"15\u00f79=1, 6\", \"31\u00f77=4, 3\"),\n  @(\"23\u00f73=7, 2\", \"82\u00f78=10, 2\"),\n  @(\"55\u00f72=27, 1\", \"18\u00f74=4, 2\"),\n  @(\"91\u00f72=45, 1\", \"79\u00f73=26, 1\"),\n  @(\"20\u00f79=2, 2\", \"74\u00f73=24, 2\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    Write-Output \"WARNING: could not find '$oldText'\"\n  }\n}\n"}
